$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "not assigned" -> "not found"
$ws.Range("B2").Value = "not found"

# Rows where "not assigned" -> "(woody)"
$woodyRows = @(4, 5, 6, 7, 52, 55, 57, 84, 147, 151, 152, 153, 160, 161, 164)
foreach ($r in $woodyRows) {
    $ws.Range("B$r").Value = "(woody)"
}
